$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110. This pushes the existing rows 110-170
# down to 111-171 (matching the diff, which shows every row from 110..170
# shifting down by one, with a brand-new row ending up at 171 that holds
# what used to be in row 170).
$ws.Rows(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44460
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112045
$ws.Cells.Item(110, 7).Value = "Zapallo"
$ws.Cells.Item(110, 8).Value = "Paine"
$ws.Cells.Item(110, 9).Value = "1a (guarda)"
$ws.Cells.Item(110, 10).Value = 750
$ws.Cells.Item(110, 11).Value = 600
$ws.Cells.Item(110, 12).Value = 600
$ws.Cells.Item(110, 13).Value = 600
$ws.Cells.Item(110, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 600
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"
